$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-shuffle the F:V (match data) content for rows whose underlying match changed position ---
$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Abha'
$arr[0,1] = 1
$arr[0,2] = 'Al Hilal'
$arr[0,3] = 3
$arr[0,4] = 5.36
$arr[0,5] = '07/08/2023 17:42'
$arr[0,6] = 5.5
$arr[0,7] = '14/08/2023 16:58'
$arr[0,8] = 4.35
$arr[0,9] = '07/08/2023 17:42'
$arr[0,10] = 4.5
$arr[0,11] = '14/08/2023 16:58'
$arr[0,12] = 1.57
$arr[0,13] = '07/08/2023 17:42'
$arr[0,14] = 1.55
$arr[0,15] = '14/08/2023 16:55'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/abha-al-hilal/xrkhHFLo/'
$ws.Range("F7:V7").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Raed'
$arr[0,1] = 0
$arr[0,2] = 'Al Ittihad'
$arr[0,3] = 3
$arr[0,4] = 4.56
$arr[0,5] = '07/08/2023 17:42'
$arr[0,6] = 15.99
$arr[0,7] = '14/08/2023 16:57'
$arr[0,8] = 4.11
$arr[0,9] = '07/08/2023 17:42'
$arr[0,10] = 8.29
$arr[0,11] = '14/08/2023 16:57'
$arr[0,12] = 1.7
$arr[0,13] = '07/08/2023 17:42'
$arr[0,14] = 1.15
$arr[0,15] = '14/08/2023 16:09'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-ittihad/jogdGZzh/'
$ws.Range("F8:V8").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ahli SC'
$arr[0,1] = 1
$arr[0,2] = 'Al Akhdoud'
$arr[0,3] = 0
$arr[0,4] = 1.16
$arr[0,5] = '22/08/2023 07:46'
$arr[0,6] = 1.11
$arr[0,7] = '24/08/2023 19:40'
$arr[0,8] = 7.91
$arr[0,9] = '22/08/2023 07:46'
$arr[0,10] = 9.28
$arr[0,11] = '24/08/2023 19:52'
$arr[0,12] = 16.84
$arr[0,13] = '22/08/2023 07:46'
$arr[0,14] = 29.41
$arr[0,15] = '24/08/2023 19:52'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ahli-sc-al-akhdoud/IgjeCGQr/'
$ws.Range("F21:V21").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Riyadh'
$arr[0,1] = 0
$arr[0,2] = 'Al Ittihad'
$arr[0,3] = 4
$arr[0,4] = 13.51
$arr[0,5] = '22/08/2023 07:46'
$arr[0,6] = 12.1
$arr[0,7] = '24/08/2023 19:54'
$arr[0,8] = 6.52
$arr[0,9] = '22/08/2023 07:46'
$arr[0,10] = 6.94
$arr[0,11] = '24/08/2023 19:54'
$arr[0,12] = 1.21
$arr[0,13] = '22/08/2023 07:46'
$arr[0,14] = 1.22
$arr[0,15] = '24/08/2023 19:07'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-riyadh-al-ittihad/zVp0Bztk/'
$ws.Range("F22:V22").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Raed'
$arr[0,1] = 0
$arr[0,2] = 'Al Hilal'
$arr[0,3] = 4
$arr[0,4] = 12.64
$arr[0,5] = '22/08/2023 07:46'
$arr[0,6] = 18.25
$arr[0,7] = '24/08/2023 19:59'
$arr[0,8] = 6.75
$arr[0,9] = '22/08/2023 07:46'
$arr[0,10] = 9.07
$arr[0,11] = '24/08/2023 19:59'
$arr[0,12] = 1.21
$arr[0,13] = '22/08/2023 07:46'
$arr[0,14] = 1.14
$arr[0,15] = '24/08/2023 19:51'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-hilal/MN4PHx3L/'
$ws.Range("F23:V23").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Nassr'
$arr[0,1] = 4
$arr[0,2] = 'Al Shabab'
$arr[0,3] = 0
$arr[0,4] = 1.39
$arr[0,5] = '27/08/2023 10:51'
$arr[0,6] = 1.35
$arr[0,7] = '29/08/2023 19:59'
$arr[0,8] = 5.14
$arr[0,9] = '27/08/2023 10:51'
$arr[0,10] = 5.69
$arr[0,11] = '29/08/2023 19:59'
$arr[0,12] = 6.68
$arr[0,13] = '27/08/2023 10:51'
$arr[0,14] = 7.34
$arr[0,15] = '29/08/2023 19:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-nassr-al-shabab/feal8PVO/'
$ws.Range("F35:V35").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ahli SC'
$arr[0,1] = 2
$arr[0,2] = 'Al Taee'
$arr[0,3] = 0
$arr[0,4] = 1.26
$arr[0,5] = '27/08/2023 10:51'
$arr[0,6] = 1.14
$arr[0,7] = '29/08/2023 19:42'
$arr[0,8] = 6.07
$arr[0,9] = '27/08/2023 10:51'
$arr[0,10] = 8.74
$arr[0,11] = '29/08/2023 19:52'
$arr[0,12] = 9.29
$arr[0,13] = '27/08/2023 10:51'
$arr[0,14] = 15.83
$arr[0,15] = '29/08/2023 19:52'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ahli-sc-al-taee/G40p95GI/'
$ws.Range("F36:V36").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Khaleej'
$arr[0,1] = 1
$arr[0,2] = 'Al Hazem'
$arr[0,3] = 1
$arr[0,4] = 2.13
$arr[0,5] = '27/08/2023 10:52'
$arr[0,6] = 2.05
$arr[0,7] = '29/08/2023 19:34'
$arr[0,8] = 3.54
$arr[0,9] = '27/08/2023 10:52'
$arr[0,10] = 3.67
$arr[0,11] = '29/08/2023 19:54'
$arr[0,12] = 3.28
$arr[0,13] = '27/08/2023 10:52'
$arr[0,14] = 3.5
$arr[0,15] = '29/08/2023 19:54'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-khaleej-al-hazem-rass/zgqmAo1C/'
$ws.Range("F37:V37").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Feiha'
$arr[0,1] = 0
$arr[0,2] = 'Al Raed'
$arr[0,3] = 0
$arr[0,4] = 1.76
$arr[0,5] = '28/08/2023 22:16'
$arr[0,6] = 2.3
$arr[0,7] = '01/09/2023 16:52'
$arr[0,8] = 4
$arr[0,9] = '28/08/2023 22:16'
$arr[0,10] = 3.37
$arr[0,11] = '01/09/2023 16:52'
$arr[0,12] = 4.13
$arr[0,13] = '28/08/2023 22:16'
$arr[0,14] = 3.2
$arr[0,15] = '01/09/2023 16:52'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-feiha-al-raed/8bmWW88P/'
$ws.Range("F38:V38").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taee'
$arr[0,1] = 1
$arr[0,2] = 'Abha'
$arr[0,3] = 0
$arr[0,4] = 1.87
$arr[0,5] = '29/08/2023 20:01'
$arr[0,6] = 2.38
$arr[0,7] = '01/09/2023 16:57'
$arr[0,8] = 4
$arr[0,9] = '29/08/2023 20:01'
$arr[0,10] = 3.24
$arr[0,11] = '01/09/2023 16:57'
$arr[0,12] = 3.66
$arr[0,13] = '29/08/2023 20:01'
$arr[0,14] = 3.18
$arr[0,15] = '01/09/2023 16:57'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-abha/pKvbSSGt/'
$ws.Range("F39:V39").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Riyadh'
$arr[0,1] = 0
$arr[0,2] = 'Al Akhdoud'
$arr[0,3] = 1
$arr[0,4] = 2.4
$arr[0,5] = '29/08/2023 17:01'
$arr[0,6] = 2.54
$arr[0,7] = '01/09/2023 19:53'
$arr[0,8] = 3.64
$arr[0,9] = '29/08/2023 17:01'
$arr[0,10] = 3.47
$arr[0,11] = '01/09/2023 19:53'
$arr[0,12] = 2.75
$arr[0,13] = '29/08/2023 17:01'
$arr[0,14] = 2.76
$arr[0,15] = '01/09/2023 19:53'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-riyadh-al-akhdoud/SrU2RnWn/'
$ws.Range("F40:V40").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ittihad'
$arr[0,1] = 3
$arr[0,2] = 'Al Hilal'
$arr[0,3] = 4
$arr[0,4] = 2.22
$arr[0,5] = '28/08/2023 22:16'
$arr[0,6] = 2.55
$arr[0,7] = '01/09/2023 19:58'
$arr[0,8] = 3.66
$arr[0,9] = '28/08/2023 22:16'
$arr[0,10] = 3.94
$arr[0,11] = '01/09/2023 19:58'
$arr[0,12] = 3.02
$arr[0,13] = '28/08/2023 22:16'
$arr[0,14] = 2.51
$arr[0,15] = '01/09/2023 19:58'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ittihad-al-hilal/6ZS6Q6og/'
$ws.Range("F41:V41").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Hazem'
$arr[0,1] = 1
$arr[0,2] = 'Al Nassr'
$arr[0,3] = 5
$arr[0,4] = 11.14
$arr[0,5] = '29/08/2023 20:01'
$arr[0,6] = 17.04
$arr[0,7] = '02/09/2023 19:56'
$arr[0,8] = 6.35
$arr[0,9] = '29/08/2023 20:01'
$arr[0,10] = 8.78
$arr[0,11] = '02/09/2023 19:56'
$arr[0,12] = 1.26
$arr[0,13] = '29/08/2023 20:01'
$arr[0,14] = 1.14
$arr[0,15] = '02/09/2023 18:32'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hazem-rass-al-nassr/f9vNMOpI/'
$ws.Range("F42:V42").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Fateh'
$arr[0,1] = 5
$arr[0,2] = 'Al Ahli SC'
$arr[0,3] = 1
$arr[0,4] = 5.41
$arr[0,5] = '29/08/2023 20:01'
$arr[0,6] = 6.86
$arr[0,7] = '02/09/2023 19:58'
$arr[0,8] = 4.35
$arr[0,9] = '29/08/2023 20:01'
$arr[0,10] = 5.34
$arr[0,11] = '02/09/2023 19:58'
$arr[0,12] = 1.6
$arr[0,13] = '29/08/2023 20:01'
$arr[0,14] = 1.39
$arr[0,15] = '02/09/2023 19:58'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-al-ahli-sc/EaWEOpG5/'
$ws.Range("F43:V43").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ettifaq'
$arr[0,1] = 3
$arr[0,2] = 'Damac'
$arr[0,3] = 1
$arr[0,4] = 1.83
$arr[0,5] = '29/08/2023 17:01'
$arr[0,6] = 1.67
$arr[0,7] = '02/09/2023 19:51'
$arr[0,8] = 3.94
$arr[0,9] = '29/08/2023 17:01'
$arr[0,10] = 3.84
$arr[0,11] = '02/09/2023 19:51'
$arr[0,12] = 4.19
$arr[0,13] = '29/08/2023 17:01'
$arr[0,14] = 5.32
$arr[0,15] = '02/09/2023 19:51'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ettifaq-fc-damac/z9XAPQ0a/'
$ws.Range("F44:V44").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Shabab'
$arr[0,1] = 1
$arr[0,2] = 'Al Khaleej'
$arr[0,3] = 3
$arr[0,4] = 1.44
$arr[0,5] = '29/08/2023 20:01'
$arr[0,6] = 1.94
$arr[0,7] = '02/09/2023 19:51'
$arr[0,8] = 4.85
$arr[0,9] = '29/08/2023 20:01'
$arr[0,10] = 3.76
$arr[0,11] = '02/09/2023 19:51'
$arr[0,12] = 7
$arr[0,13] = '29/08/2023 20:01'
$arr[0,14] = 3.78
$arr[0,15] = '02/09/2023 19:51'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-shabab-al-khaleej/dOLJN4VB/'
$ws.Range("F45:V45").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taawon'
$arr[0,1] = 4
$arr[0,2] = 'Al Wehda'
$arr[0,3] = 1
$arr[0,4] = 1.83
$arr[0,5] = '29/08/2023 17:01'
$arr[0,6] = 2.37
$arr[0,7] = '02/09/2023 19:59'
$arr[0,8] = 3.94
$arr[0,9] = '29/08/2023 17:01'
$arr[0,10] = 3.48
$arr[0,11] = '02/09/2023 19:59'
$arr[0,12] = 4.19
$arr[0,13] = '29/08/2023 17:01'
$arr[0,14] = 2.98
$arr[0,15] = '02/09/2023 19:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taawon-al-wehda/lWwRLr0O/'
$ws.Range("F46:V46").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ettifaq'
$arr[0,1] = 4
$arr[0,2] = 'Al Taee'
$arr[0,3] = 3
$arr[0,4] = 1.5
$arr[0,5] = '19/09/2023 19:36'
$arr[0,6] = 1.47
$arr[0,7] = '21/09/2023 16:56'
$arr[0,8] = 4.15
$arr[0,9] = '19/09/2023 19:36'
$arr[0,10] = 4.55
$arr[0,11] = '21/09/2023 16:56'
$arr[0,12] = 5.75
$arr[0,13] = '19/09/2023 19:36'
$arr[0,14] = 6.56
$arr[0,15] = '21/09/2023 16:56'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ettifaq-fc-al-taee/jaMRj3ab/'
$ws.Range("F57:V57").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taawon'
$arr[0,1] = 2
$arr[0,2] = 'Al Raed'
$arr[0,3] = 1
$arr[0,4] = 1.57
$arr[0,5] = '19/09/2023 19:36'
$arr[0,6] = 2.01
$arr[0,7] = '21/09/2023 16:51'
$arr[0,8] = 4.12
$arr[0,9] = '19/09/2023 19:36'
$arr[0,10] = 3.74
$arr[0,11] = '21/09/2023 16:51'
$arr[0,12] = 4.95
$arr[0,13] = '19/09/2023 19:36'
$arr[0,14] = 3.57
$arr[0,15] = '21/09/2023 16:51'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taawon-al-raed/tUvyUuho/'
$ws.Range("F58:V58").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Shabab'
$arr[0,1] = 4
$arr[0,2] = 'Al Hazem'
$arr[0,3] = 1
$arr[0,4] = 1.65
$arr[0,5] = '19/09/2023 19:38'
$arr[0,6] = 1.44
$arr[0,7] = '21/09/2023 19:53'
$arr[0,8] = 3.83
$arr[0,9] = '19/09/2023 19:38'
$arr[0,10] = 4.92
$arr[0,11] = '21/09/2023 19:53'
$arr[0,12] = 4.72
$arr[0,13] = '19/09/2023 19:38'
$arr[0,14] = 6.57
$arr[0,15] = '21/09/2023 19:53'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-shabab-al-hazem-rass/fXSIhPUo/'
$ws.Range("F59:V59").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ittihad'
$arr[0,1] = 2
$arr[0,2] = 'Al Fateh'
$arr[0,3] = 1
$arr[0,4] = 1.23
$arr[0,5] = '19/09/2023 19:37'
$arr[0,6] = 1.71
$arr[0,7] = '21/09/2023 19:54'
$arr[0,8] = 6.23
$arr[0,9] = '19/09/2023 19:37'
$arr[0,10] = 4.05
$arr[0,11] = '21/09/2023 19:55'
$arr[0,12] = 8.61
$arr[0,13] = '19/09/2023 19:37'
$arr[0,14] = 4.57
$arr[0,15] = '21/09/2023 19:55'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ittihad-al-fateh/YNRMiqph/'
$ws.Range("F60:V60").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ahli SC'
$arr[0,1] = 3
$arr[0,2] = 'Al Wehda'
$arr[0,3] = 1
$arr[0,4] = 1.22
$arr[0,5] = '14/10/2023 19:16'
$arr[0,6] = 1.39
$arr[0,7] = '21/10/2023 19:18'
$arr[0,8] = 6.35
$arr[0,9] = '14/10/2023 19:16'
$arr[0,10] = 5.04
$arr[0,11] = '21/10/2023 19:37'
$arr[0,12] = 9.09
$arr[0,13] = '14/10/2023 19:16'
$arr[0,14] = 7.62
$arr[0,15] = '21/10/2023 19:37'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ahli-sc-al-wehda/CUpZxdhD/'
$ws.Range("F89:V89").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Shabab'
$arr[0,1] = 2
$arr[0,2] = 'Al Taee'
$arr[0,3] = 0
$arr[0,4] = 1.5
$arr[0,5] = '14/10/2023 19:16'
$arr[0,6] = 1.52
$arr[0,7] = '21/10/2023 19:59'
$arr[0,8] = 4.49
$arr[0,9] = '14/10/2023 19:16'
$arr[0,10] = 4.52
$arr[0,11] = '21/10/2023 19:59'
$arr[0,12] = 6.2
$arr[0,13] = '14/10/2023 19:16'
$arr[0,14] = 5.78
$arr[0,15] = '21/10/2023 19:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-shabab-al-taee/2a9iew8s/'
$ws.Range("F90:V90").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Raed'
$arr[0,1] = 1
$arr[0,2] = 'Al Fateh'
$arr[0,3] = 2
$arr[0,4] = 4.14
$arr[0,5] = '24/10/2023 22:01'
$arr[0,6] = 2.93
$arr[0,7] = '27/10/2023 16:59'
$arr[0,8] = 4.26
$arr[0,9] = '24/10/2023 22:01'
$arr[0,10] = 3.98
$arr[0,11] = '27/10/2023 16:59'
$arr[0,12] = 1.75
$arr[0,13] = '24/10/2023 22:01'
$arr[0,14] = 2.21
$arr[0,15] = '27/10/2023 16:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-fateh/lGeDjFMC/'
$ws.Range("F95:V95").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taee'
$arr[0,1] = 3
$arr[0,2] = 'Al Riyadh'
$arr[0,3] = 2
$arr[0,4] = 1.84
$arr[0,5] = '24/10/2023 22:01'
$arr[0,6] = 1.75
$arr[0,7] = '27/10/2023 16:59'
$arr[0,8] = 3.78
$arr[0,9] = '24/10/2023 22:01'
$arr[0,10] = 3.85
$arr[0,11] = '27/10/2023 16:59'
$arr[0,12] = 4.13
$arr[0,13] = '24/10/2023 22:01'
$arr[0,14] = 4.65
$arr[0,15] = '27/10/2023 16:58'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-riyadh/pxvV8dxs/'
$ws.Range("F96:V96").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taee'
$arr[0,1] = 3
$arr[0,2] = 'Al Feiha'
$arr[0,3] = 3
$arr[0,4] = 2.73
$arr[0,5] = '29/10/2023 19:43'
$arr[0,6] = 3.06
$arr[0,7] = '03/11/2023 15:55'
$arr[0,8] = 3.38
$arr[0,9] = '29/10/2023 19:43'
$arr[0,10] = 3.43
$arr[0,11] = '03/11/2023 15:55'
$arr[0,12] = 2.45
$arr[0,13] = '29/10/2023 19:43'
$arr[0,14] = 2.35
$arr[0,15] = '03/11/2023 15:55'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-feiha/SKZk4YjC/'
$ws.Range("F101:V101").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Fateh'
$arr[0,1] = 0
$arr[0,2] = 'Al Hilal'
$arr[0,3] = 2
$arr[0,4] = 5.32
$arr[0,5] = '29/10/2023 19:43'
$arr[0,6] = 10.37
$arr[0,7] = '03/11/2023 15:57'
$arr[0,8] = 5.18
$arr[0,9] = '29/10/2023 19:43'
$arr[0,10] = 6.9
$arr[0,11] = '03/11/2023 15:57'
$arr[0,12] = 1.43
$arr[0,13] = '29/10/2023 19:43'
$arr[0,14] = 1.23
$arr[0,15] = '03/11/2023 15:52'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-al-hilal/0YXc2CLO/'
$ws.Range("F102:V102").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Akhdoud'
$arr[0,1] = 2
$arr[0,2] = 'Al Hazem'
$arr[0,3] = 1
$arr[0,4] = 2.19
$arr[0,5] = '08/11/2023 06:12'
$arr[0,6] = 1.72
$arr[0,7] = '11/11/2023 15:56'
$arr[0,8] = 3.47
$arr[0,9] = '08/11/2023 06:12'
$arr[0,10] = 4.09
$arr[0,11] = '11/11/2023 15:56'
$arr[0,12] = 3.06
$arr[0,13] = '08/11/2023 06:12'
$arr[0,14] = 4.49
$arr[0,15] = '11/11/2023 15:56'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-akhdoud-al-hazem-rass/AaT8Kgyg/'
$ws.Range("F116:V116").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Feiha'
$arr[0,1] = 0
$arr[0,2] = 'Al Ettifaq'
$arr[0,3] = 0
$arr[0,4] = 2.7
$arr[0,5] = '04/11/2023 19:13'
$arr[0,6] = 3.12
$arr[0,7] = '11/11/2023 15:56'
$arr[0,8] = 3.21
$arr[0,9] = '04/11/2023 19:13'
$arr[0,10] = 3.33
$arr[0,11] = '11/11/2023 15:56'
$arr[0,12] = 2.57
$arr[0,13] = '04/11/2023 19:13'
$arr[0,14] = 2.36
$arr[0,15] = '11/11/2023 15:56'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-feiha-al-ettifaq-fc/ddP4LZLn/'
$ws.Range("F117:V117").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taawon'
$arr[0,1] = 1
$arr[0,2] = 'Al Riyadh'
$arr[0,3] = 2
$arr[0,4] = 1.35
$arr[0,5] = '17/11/2023 16:42'
$arr[0,6] = 1.37
$arr[0,7] = '24/11/2023 15:50'
$arr[0,8] = 5.11
$arr[0,9] = '17/11/2023 16:42'
$arr[0,10] = 5.32
$arr[0,11] = '24/11/2023 15:54'
$arr[0,12] = 6.86
$arr[0,13] = '17/11/2023 16:42'
$arr[0,14] = 7.56
$arr[0,15] = '24/11/2023 15:54'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taawon-al-riyadh/SbM0IOBM/'
$ws.Range("F119:V119").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ettifaq'
$arr[0,1] = 1
$arr[0,2] = 'Al Ittihad'
$arr[0,3] = 1
$arr[0,4] = 3.84
$arr[0,5] = '17/11/2023 16:42'
$arr[0,6] = 3.08
$arr[0,7] = '24/11/2023 15:50'
$arr[0,8] = 3.87
$arr[0,9] = '17/11/2023 16:42'
$arr[0,10] = 3.34
$arr[0,11] = '24/11/2023 15:50'
$arr[0,12] = 1.79
$arr[0,13] = '17/11/2023 16:42'
$arr[0,14] = 2.38
$arr[0,15] = '24/11/2023 15:50'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ettifaq-fc-al-ittihad/KYIhKpsA/'
$ws.Range("F120:V120").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taee'
$arr[0,1] = 4
$arr[0,2] = 'Al Raed'
$arr[0,3] = 3
$arr[0,4] = 2.75
$arr[0,5] = '17/11/2023 19:43'
$arr[0,6] = 4.15
$arr[0,7] = '24/11/2023 18:57'
$arr[0,8] = 3.67
$arr[0,9] = '17/11/2023 19:43'
$arr[0,10] = 3.63
$arr[0,11] = '24/11/2023 18:58'
$arr[0,12] = 2.3
$arr[0,13] = '17/11/2023 19:43'
$arr[0,14] = 1.88
$arr[0,15] = '24/11/2023 18:53'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-raed/trKlLQR3/'
$ws.Range("F121:V121").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Nassr'
$arr[0,1] = 3
$arr[0,2] = 'Al Akhdoud'
$arr[0,3] = 0
$arr[0,4] = 1.12
$arr[0,5] = '17/11/2023 19:43'
$arr[0,6] = 1.24
$arr[0,7] = '24/11/2023 18:55'
$arr[0,8] = 9.47
$arr[0,9] = '17/11/2023 19:43'
$arr[0,10] = 7.01
$arr[0,11] = '24/11/2023 18:57'
$arr[0,12] = 12.33
$arr[0,13] = '17/11/2023 19:43'
$arr[0,14] = 9.31
$arr[0,15] = '24/11/2023 18:57'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-nassr-al-akhdoud/pbBvBt4d/'
$ws.Range("F122:V122").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Ahli SC'
$arr[0,1] = 0
$arr[0,2] = 'Al Shabab'
$arr[0,3] = 0
$arr[0,4] = 1.51
$arr[0,5] = '18/11/2023 19:43'
$arr[0,6] = 2.06
$arr[0,7] = '25/11/2023 18:59'
$arr[0,8] = 4.54
$arr[0,9] = '18/11/2023 19:43'
$arr[0,10] = 3.84
$arr[0,11] = '25/11/2023 18:59'
$arr[0,12] = 5.05
$arr[0,13] = '18/11/2023 19:43'
$arr[0,14] = 3.33
$arr[0,15] = '25/11/2023 18:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ahli-sc-al-shabab/67CzCMlj/'
$ws.Range("F124:V124").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Fateh'
$arr[0,1] = 0
$arr[0,2] = 'Al Feiha'
$arr[0,3] = 1
$arr[0,4] = 1.97
$arr[0,5] = '18/11/2023 18:13'
$arr[0,6] = 1.94
$arr[0,7] = '25/11/2023 18:34'
$arr[0,8] = 3.99
$arr[0,9] = '18/11/2023 18:13'
$arr[0,10] = 4.11
$arr[0,11] = '25/11/2023 18:34'
$arr[0,12] = 3.17
$arr[0,13] = '18/11/2023 18:13'
$arr[0,14] = 3.48
$arr[0,15] = '25/11/2023 18:17'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-al-feiha/ANB5HrRS/'
$ws.Range("F125:V125").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Wehda'
$arr[0,1] = 3
$arr[0,2] = 'Al Khaleej'
$arr[0,3] = 1
$arr[0,4] = 1.9
$arr[0,5] = '18/11/2023 19:43'
$arr[0,6] = 2.05
$arr[0,7] = '25/11/2023 18:59'
$arr[0,8] = 3.77
$arr[0,9] = '18/11/2023 19:43'
$arr[0,10] = 3.72
$arr[0,11] = '25/11/2023 18:59'
$arr[0,12] = 3.53
$arr[0,13] = '18/11/2023 19:43'
$arr[0,14] = 3.47
$arr[0,15] = '25/11/2023 18:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-wehda-al-khaleej/O6PqM6Cc/'
$ws.Range("F126:V126").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Damac'
$arr[0,1] = 4
$arr[0,2] = 'Abha'
$arr[0,3] = 2
$arr[0,4] = 1.59
$arr[0,5] = '18/11/2023 19:43'
$arr[0,6] = 1.9
$arr[0,7] = '25/11/2023 18:59'
$arr[0,8] = 4.17
$arr[0,9] = '18/11/2023 19:43'
$arr[0,10] = 3.9
$arr[0,11] = '25/11/2023 18:59'
$arr[0,12] = 4.75
$arr[0,13] = '18/11/2023 19:43'
$arr[0,14] = 3.78
$arr[0,15] = '25/11/2023 18:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/damac-abha/EeIdJ4dG/'
$ws.Range("F127:V127").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Akhdoud'
$arr[0,1] = 1
$arr[0,2] = 'Al Ettifaq'
$arr[0,3] = 0
$arr[0,4] = 3.25
$arr[0,5] = '28/11/2023 11:42'
$arr[0,6] = 3.2
$arr[0,7] = '02/12/2023 15:57'
$arr[0,8] = 3.59
$arr[0,9] = '28/11/2023 11:42'
$arr[0,10] = 3.4
$arr[0,11] = '02/12/2023 15:57'
$arr[0,12] = 2.06
$arr[0,13] = '28/11/2023 11:42'
$arr[0,14] = 2.28
$arr[0,15] = '02/12/2023 15:57'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-akhdoud-al-ettifaq-fc/06gRYqd3/'
$ws.Range("F133:V133").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Fateh'
$arr[0,1] = 0
$arr[0,2] = 'Al Taee'
$arr[0,3] = 1
$arr[0,4] = 1.61
$arr[0,5] = '25/11/2023 19:13'
$arr[0,6] = 1.58
$arr[0,7] = '02/12/2023 15:56'
$arr[0,8] = 4.31
$arr[0,9] = '25/11/2023 19:13'
$arr[0,10] = 4.57
$arr[0,11] = '02/12/2023 15:56'
$arr[0,12] = 5.01
$arr[0,13] = '25/11/2023 19:13'
$arr[0,14] = 5.02
$arr[0,15] = '02/12/2023 15:56'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-al-taee/n3kNZPtc/'
$ws.Range("F134:V134").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Raed'
$arr[0,1] = 2
$arr[0,2] = 'Al Wehda'
$arr[0,3] = 0
$arr[0,4] = 2.39
$arr[0,5] = '25/11/2023 21:43'
$arr[0,6] = 3.04
$arr[0,7] = '02/12/2023 18:58'
$arr[0,8] = 3.41
$arr[0,9] = '25/11/2023 21:43'
$arr[0,10] = 3.52
$arr[0,11] = '02/12/2023 18:58'
$arr[0,12] = 2.79
$arr[0,13] = '25/11/2023 21:43'
$arr[0,14] = 2.32
$arr[0,15] = '02/12/2023 18:58'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-wehda/jenvWstM/'
$ws.Range("F135:V135").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Shabab'
$arr[0,1] = 1
$arr[0,2] = 'Al Taawon'
$arr[0,3] = 2
$arr[0,4] = 2.17
$arr[0,5] = '25/11/2023 21:43'
$arr[0,6] = 2.23
$arr[0,7] = '02/12/2023 18:18'
$arr[0,8] = 3.41
$arr[0,9] = '25/11/2023 21:43'
$arr[0,10] = 3.45
$arr[0,11] = '02/12/2023 18:59'
$arr[0,12] = 3.15
$arr[0,13] = '25/11/2023 21:43'
$arr[0,14] = 3.26
$arr[0,15] = '02/12/2023 18:59'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-shabab-al-taawon/b1yEPuJk/'
$ws.Range("F136:V136").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taee'
$arr[0,1] = 1
$arr[0,2] = 'Al Hilal'
$arr[0,3] = 2
$arr[0,4] = 12.91
$arr[0,5] = '03/12/2023 18:12'
$arr[0,6] = 20.78
$arr[0,7] = '08/12/2023 15:59'
$arr[0,8] = 9.43
$arr[0,9] = '03/12/2023 18:12'
$arr[0,10] = 11.18
$arr[0,11] = '08/12/2023 15:59'
$arr[0,12] = 1.1
$arr[0,13] = '03/12/2023 18:12'
$arr[0,14] = 1.1
$arr[0,15] = '08/12/2023 15:25'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-hilal/4dN89pzL/'
$ws.Range("F141:V141").Value2 = $arr

$arr = New-Object "object[,]" 1,17
$arr[0,0] = 'Al Taawon'
$arr[0,1] = 4
$arr[0,2] = 'Al Feiha'
$arr[0,3] = 1
$arr[0,4] = 1.81
$arr[0,5] = '03/12/2023 18:12'
$arr[0,6] = 1.76
$arr[0,7] = '08/12/2023 15:51'
$arr[0,8] = 3.7
$arr[0,9] = '03/12/2023 18:12'
$arr[0,10] = 3.8
$arr[0,11] = '08/12/2023 15:51'
$arr[0,12] = 3.91
$arr[0,13] = '03/12/2023 18:12'
$arr[0,14] = 4.59
$arr[0,15] = '08/12/2023 15:51'
$arr[0,16] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taawon-al-feiha/Ac6i1Ms1/'
$ws.Range("F142:V142").Value2 = $arr

# --- Append new rows 157-159 (copy formatting from row 156 first) ---
$srcFmt = $ws.Range("A156:V156")
$dstFmt = $ws.Range("A157:V159")
$srcFmt.Copy() | Out-Null
$dstFmt.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$arr = New-Object "object[,]" 1,22
$arr[0,0] = 156
$arr[0,1] = 'saudi-arabia'
$arr[0,2] = 'saudi-professional-league'
$arr[0,3] = '2023-2024'
$arr[0,4] = 45282.66666666666
$arr[0,5] = 'Al Akhdoud'
$arr[0,6] = 1
$arr[0,7] = 'Al Shabab'
$arr[0,8] = 0
$arr[0,9] = 2.51
$arr[0,10] = '20/12/2023 03:12'
$arr[0,11] = 2.71
$arr[0,12] = '22/12/2023 15:56'
$arr[0,13] = 3.38
$arr[0,14] = '20/12/2023 03:12'
$arr[0,15] = 3.59
$arr[0,16] = '22/12/2023 15:25'
$arr[0,17] = 2.66
$arr[0,18] = '20/12/2023 03:12'
$arr[0,19] = 2.42
$arr[0,20] = '22/12/2023 15:25'
$arr[0,21] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-akhdoud-al-shabab/E3CtKzlt/'
$ws.Range("A157:V157").Value2 = $arr

$arr = New-Object "object[,]" 1,22
$arr[0,0] = 157
$arr[0,1] = 'saudi-arabia'
$arr[0,2] = 'saudi-professional-league'
$arr[0,3] = '2023-2024'
$arr[0,4] = 45282.66666666666
$arr[0,5] = 'Al Hazem'
$arr[0,6] = 0
$arr[0,7] = 'Al Ahli SC'
$arr[0,8] = 4
$arr[0,9] = 6.42
$arr[0,10] = '15/12/2023 17:43'
$arr[0,11] = 7.66
$arr[0,12] = '22/12/2023 15:58'
$arr[0,13] = 4.84
$arr[0,14] = '15/12/2023 17:43'
$arr[0,15] = 5.44
$arr[0,16] = '22/12/2023 15:58'
$arr[0,17] = 1.39
$arr[0,18] = '15/12/2023 17:43'
$arr[0,19] = 1.36
$arr[0,20] = '22/12/2023 15:51'
$arr[0,21] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hazem-rass-al-ahli-sc/Ie9cGhm5/'
$ws.Range("A158:V158").Value2 = $arr

$arr = New-Object "object[,]" 1,22
$arr[0,0] = 158
$arr[0,1] = 'saudi-arabia'
$arr[0,2] = 'saudi-professional-league'
$arr[0,3] = '2023-2024'
$arr[0,4] = 45282.66666666666
$arr[0,5] = 'Al Nassr'
$arr[0,6] = 3
$arr[0,7] = 'Al Ettifaq'
$arr[0,8] = 1
$arr[0,9] = 1.36
$arr[0,10] = '16/12/2023 19:13'
$arr[0,11] = 1.32
$arr[0,12] = '22/12/2023 15:51'
$arr[0,13] = 5.14
$arr[0,14] = '16/12/2023 19:13'
$arr[0,15] = 5.67
$arr[0,16] = '22/12/2023 15:56'
$arr[0,17] = 6.48
$arr[0,18] = '16/12/2023 19:13'
$arr[0,19] = 8.5
$arr[0,20] = '22/12/2023 15:56'
$arr[0,21] = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-nassr-al-ettifaq-fc/lM4gHYXb/'
$ws.Range("A159:V159").Value2 = $arr
